$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colMap = @{ "B" = 2; "C" = 3; "D" = 4; "E" = 5; "F" = 6; "K" = 11; "M" = 13; "O" = 15 }

$data = @{
    2 = @{ "B" = 8.658781320990395; "C" = 5.776349252250848; "D" = 5.210700280086404; "E" = 12.80913531481472; "F" = 25.8505598011867; "K" = 7.881343650347751; "M" = 13.37592750253274; "O" = 23.12837045790562 }
    3 = @{ "B" = 8.375255007433136; "C" = 5.677267269457328; "D" = 5.16834100440316; "E" = 12.59790147270134; "F" = 25.85822056369273; "K" = 7.690156342032411; "M" = 13.20665415481002; "O" = 23.18463731653462 }
    4 = @{ "B" = 8.197753954047924; "C" = 5.614898613606682; "D" = 5.141725974095158; "E" = 12.47089036794688; "F" = 25.86990825973226; "K" = 7.571603531575046; "M" = 13.10492469694334; "O" = 23.2243900530789 }
    5 = @{ "B" = 8.124679888736717; "C" = 5.589116164459362; "D" = 5.130731789761617; "E" = 12.41988153139816; "F" = 25.87642538715141; "K" = 7.52307155495469; "M" = 13.06407271075935; "O" = 23.2418948995931 }
    6 = @{ "B" = 8.112504868418501; "C" = 5.584813414568698; "D" = 5.128897375259167; "E" = 12.41145892354472; "F" = 25.8776134301305; "K" = 7.515001680833103; "M" = 13.05732713711795; "O" = 23.24488029575614 }
    7 = @{ "B" = 8.196771289259223; "C" = 5.614552362320814; "D" = 5.14157829786653; "E" = 12.47019931532193; "F" = 25.86998905223733; "K" = 7.570949809122407; "M" = 13.10437124483402; "O" = 23.22462084931468 }
    8 = @{ "B" = 8.56179788794508; "C" = 5.742514597092134; "D" = 5.196223074784131; "E" = 12.73578857916091; "F" = 25.85175090478663; "K" = 7.81570568391245; "M" = 13.31713419112207; "O" = 23.14668902793693 }
    9 = @{ "B" = 9.245829622695178; "C" = 5.980520859281305; "D" = 5.298395032029184; "E" = 13.27465665702039; "F" = 25.8714380582451; "K" = 8.36248194742536; "M" = 13.74965334871002; "O" = 23.03530757062607 }
    10 = @{ "B" = 9.723405655031021; "C" = 6.146528865376229; "D" = 5.370192279400342; "E" = 13.67717495384462; "F" = 25.91969059850083; "K" = 8.810338483946493; "M" = 14.07375893708164; "O" = 22.97892586175816 }
    11 = @{ "B" = 9.934245163249116; "C" = 6.219938267475892; "D" = 5.402097980420135; "E" = 13.86081283664143; "F" = 25.94895079484193; "K" = 9.005148249071199; "M" = 14.22194387869306; "O" = 22.95883672461174 }
    12 = @{ "B" = 10.013091296946; "C" = 6.24741862179135; "D" = 5.414067340635546; "E" = 13.93035423698366; "F" = 25.96107758436969; "K" = 9.077603665385309; "M" = 14.27811290257574; "O" = 22.95203116339198 }
    13 = @{ "B" = 9.996155650796096; "C" = 6.241514616605433; "D" = 5.411494597515208; "E" = 13.91537833504096; "F" = 25.95841940169086; "K" = 9.062058029598006; "M" = 14.2660143010042; "O" = 22.95346117032844 }
    14 = @{ "B" = 9.940752191662751; "C" = 6.222205552481937; "D" = 5.403084987063426; "E" = 13.86653443284874; "F" = 25.94992752843156; "K" = 9.011135716467093; "M" = 14.226564143105; "O" = 22.95826074623421 }
    15 = @{ "B" = 9.906684496817618; "C" = 6.210336325532372; "D" = 5.397919071127089; "E" = 13.83661418771482; "F" = 25.94486215114597; "K" = 8.979772236542779; "M" = 14.20240534055922; "O" = 22.96130510202127 }
    16 = @{ "B" = 9.709491410553435; "C" = 6.141687657596496; "D" = 5.36809162256842; "E" = 13.66517802889427; "F" = 25.91792510725512; "K" = 8.797424843315092; "M" = 14.06408520191277; "O" = 22.9803508079257 }
    17 = @{ "B" = 9.586823531943194; "C" = 6.099022932730411; "D" = 5.349597116111712; "E" = 13.56008951441965; "F" = 25.90326948173583; "K" = 8.68325127934518; "M" = 13.97938382242237; "O" = 22.9934603515502 }
    18 = @{ "B" = 9.515668429300217; "C" = 6.074285690959334; "D" = 5.338888715160336; "E" = 13.49969961181533; "F" = 25.89552834469653; "K" = 8.616743432059074; "M" = 13.93073931780893; "O" = 23.00152375059898 }
    19 = @{ "B" = 9.4914758554531; "C" = 6.065876616346193; "D" = 5.335250990556491; "E" = 13.47926425408729; "F" = 25.89302567945679; "K" = 8.594082095037065; "M" = 13.91428337447593; "O" = 23.00434365105811 }
    20 = @{ "B" = 9.599944361180023; "C" = 6.103585230038664; "D" = 5.351573250129578; "E" = 13.57127130132179; "F" = 25.90475838130905; "K" = 8.695492248921731; "M" = 13.98839320958076; "O" = 22.99201065747781 }
    21 = @{ "B" = 9.957053047487586; "C" = 6.227885843333144; "D" = 5.405558182107285; "E" = 13.88088160554092; "F" = 25.95239343580305; "K" = 9.026128748589917; "M" = 14.23815054532703; "O" = 22.95682921780335 }
    22 = @{ "B" = 10.18462139694663; "C" = 6.307262494278858; "D" = 5.440182122186831; "E" = 14.08320718466635; "F" = 25.98962303251533; "K" = 9.234543444498209; "M" = 14.40167592904987; "O" = 22.93851023583612 }
    23 = @{ "B" = 10.06371841312131; "C" = 6.265072680778516; "D" = 5.421764219497752; "E" = 13.9752482300922; "F" = 25.96919681400886; "K" = 9.12402024052721; "M" = 14.31438959959531; "O" = 22.94785904648485 }
    24 = @{ "B" = 9.59401439633954; "C" = 6.101523264426769; "D" = 5.350680075295744; "E" = 13.56621592603092; "F" = 25.90408311703188; "K" = 8.689960807615828; "M" = 13.98431990134797; "O" = 22.99266442414358 }
    25 = @{ "B" = 9.064799265140158; "C" = 5.917621717070708; "D" = 5.271312511361257; "E" = 13.12740206836699; "F" = 25.86017267072875; "K" = 8.189881429539168; "M" = 13.63131638324489; "O" = 23.06098219567858 }
}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $colIdx = $colMap[$c]
        $ws.Cells.Item([int]$r, $colIdx).Value = $data[$r][$c]
    }
}
